$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.0001998
$ws.Range("F2").Value = 0.01541187
$ws.Range("G2").Value = 0.000375365664

$ws.Range("E3").Value = 0.00415737
$ws.Range("F3").Value = 0.00930933
$ws.Range("G3").Value = 0.004853213702623906
